$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update data values ---
# Row 2 (id=1, ABOBAKAR): email + message text updated
$ws.Range("C2").Value = "ranaabobakarit@gmail.com"
$ws.Range("D2").Value = "kya hal ha"

# Row 3 (id=2, SHAHZEB): email unchanged value, message updated
$ws.Range("C3").Value = "ranaabobakar777@gmail.com"
$ws.Range("D3").Value = "kya hal ha"

# Row 4 (id=3, ZEESHAN): email + message updated
$ws.Range("C4").Value = "abobakarit786@gmail.com"
$ws.Range("D4").Value = "kya hal ha"

# --- Apply new formatting to D3 and D4 (new font + wrap/shrink alignment) ---
$ws.Range("D3").Font.Name = "Arial"
$ws.Range("D3").Font.Size = 10
$ws.Range("D3").VerticalAlignment = -4107
$ws.Range("D3").WrapText = $true
$ws.Range("D3").ShrinkToFit = $true

$ws.Range("D4").Font.Name = "Arial"
$ws.Range("D4").Font.Size = 10
$ws.Range("D4").VerticalAlignment = -4107
$ws.Range("D4").WrapText = $true
$ws.Range("D4").ShrinkToFit = $true

# --- Remove old last row (id=4, shan) entirely ---
$ws.Rows.Item(5).Delete()

# --- Widen column D ---
$ws.Columns.Item(4).ColumnWidth = 34.7

# --- Update selection to match saved view state ---
$ws.Range("D8").Select()
